# "Generate Report for Handoff"
# b.md has moved from "Handed back: in sync with en-US" to "Ready for handoff",
# with a freshly generated handoff package (b.63290e5768f688058c7b37413b0a5c26c308f864.*.xlf)
# and a new handoff timestamp. Update the Overview sheet and both per-locale sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: row for b.md (row 3)
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value2 = "Ready for handoff"
$overview.Range("C3").Value2 = "Ready for handoff"
$overview.Range("D3").Value2 = "2016-24-18 12:24:15"

# ---------------------------------------------------------------------------
# zh-cn sheet: row for b.md (row 3)
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value2 = "Ready for handoff"
$zhcn.Range("D3").Value2 = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value2 = "2016-03-18 12:24:12"

# Rebuild the hyperlinks on the zh-cn sheet so the D3 hyperlink's display text
# matches the new target filename (the underlying link target is unchanged).
$zhcnA2 = "https://github.com/OpenLocalizationTest/oltest/blob/108e1090d0794f6bc17451c43364088dd1758322/e2e/a.md"
$zhcnD2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/422b50eb09acb845a8102fbca7ada86b564c7a72/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcnF2 = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/3fb01d93ee5bff617b8c18c72ba054fb22cf799f/e2e/a.md"
$zhcnG2 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d29690c13e9a38fd0ed3d28ac3c8dc6293900415/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcnA3 = "https://github.com/OpenLocalizationTest/oltest/blob/108e1090d0794f6bc17451c43364088dd1758322/e2e/b.md"

$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), $zhcnA2, "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("B2"), $zhcnA2, "", "", ".md")
$zhcn.Hyperlinks.Add($zhcn.Range("D2"), $zhcnD2, "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), $zhcnF2, "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), $zhcnG2, "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), $zhcnA3, "", "", "b.md")
$zhcn.Hyperlinks.Add($zhcn.Range("B3"), $zhcnA3, "", "", ".md")
$zhcn.Hyperlinks.Add($zhcn.Range("D3"), $zhcnD2, "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), $zhcnF2, "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), $zhcnG2, "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")

# ---------------------------------------------------------------------------
# de-de sheet: row for b.md (row 3)
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value2 = "Ready for handoff"
$dede.Range("D3").Value2 = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value2 = "2016-03-18 12:24:15"

# Rebuild the hyperlinks on the de-de sheet the same way.
$dedeA2 = "https://github.com/OpenLocalizationTest/oltest/blob/108e1090d0794f6bc17451c43364088dd1758322/e2e/a.md"
$dedeD2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0b9f80ae6875ebb883bb2f5ba85e1e75dcb73e27/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dedeF2 = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/994afd6b27c2d1d302bf9aea04a51c62392ddfbe/e2e/a.md"
$dedeG2 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/124a6ea631102a755acf1c9bf68dfc87c60dad51/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dedeA3 = "https://github.com/OpenLocalizationTest/oltest/blob/108e1090d0794f6bc17451c43364088dd1758322/e2e/b.md"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), $dedeA2, "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("B2"), $dedeA2, "", "", ".md")
$dede.Hyperlinks.Add($dede.Range("D2"), $dedeD2, "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("F2"), $dedeF2, "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("G2"), $dedeG2, "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("A3"), $dedeA3, "", "", "b.md")
$dede.Hyperlinks.Add($dede.Range("B3"), $dedeA3, "", "", ".md")
$dede.Hyperlinks.Add($dede.Range("D3"), $dedeD2, "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("F3"), $dedeF2, "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("G3"), $dedeG2, "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
